$d = $word.ActiveDocument

# Row 1 of the first table holds the "Month" / "Savings" header cells.
# The old Apache POI 4.1.0 writer emitted the XML Boolean values
# true/false for the w:b / w:i / w:strike run-property toggles; POI
# 5.2.3 emits on/off tokens instead. Re-apply the same (semantically
# unchanged) bold / not-italic / not-struck-through formatting to those
# two header runs so they get re-serialized with the current writer's
# tokens.

$table = $d.Tables.Item(1)

for ($col = 1; $col -le 2; $col++) {
    $cell = $table.Cell(1, $col)
    $cellRange = $cell.Range

    # Exclude the trailing end-of-cell mark so only the header text
    # run's character formatting is touched (not the cell mark itself).
    $cellRange.MoveEnd(1, -1) | Out-Null

    $cellRange.Font.Bold = $true
    $cellRange.Font.Italic = $false
    $cellRange.Font.StrikeThrough = $false
}
